$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.872027
$ws.Range("H2").Value = 59.61608099999999
$ws.Range("I2").Value = 0.2756064822985579
$ws.Range("J2").Value = 0.2756064822985579
$ws.Range("M2").Value = 4.224096333333333
$ws.Range("N2").Value = 12.672289
$ws.Range("O2").Value = 0.1277189908446358
$ws.Range("P2").Value = 0.1277189908446358
$ws.Range("Q2").Value = 83.94135638660099
$ws.Range("R2").Value = 755.4722074794089
$ws.Range("S2").Value = 0.0352001817894118
$ws.Range("T2").Value = 0.0352001817894118
$ws.Range("G3").Value = 19.872027
$ws.Range("H3").Value = 59.61608099999999
$ws.Range("I3").Value = 0.2756064822985579
$ws.Range("J3").Value = 0.2756064822985579
$ws.Range("O3").Value = 0.4492078640046304
$ws.Range("P3").Value = 0.4492078640046304
$ws.Range("Q3").Value = 295.235008941978
$ws.Range("R3").Value = 2657.115080477802
$ws.Range("S3").Value = 0.1238045992191652
$ws.Range("T3").Value = 0.1238045992191652
$ws.Range("G4").Value = 19.872027
$ws.Range("H4").Value = 59.61608099999999
$ws.Range("I4").Value = 0.2756064822985579
$ws.Range("J4").Value = 0.2756064822985579
$ws.Range("O4").Value = 0.4230731451507339
$ws.Range("P4").Value = 0.4230731451507338
$ws.Range("Q4").Value = 278.058364068177
$ws.Range("R4").Value = 2502.525276613593
$ws.Range("S4").Value = 0.1166017012899809
$ws.Range("T4").Value = 0.1166017012899809
$ws.Range("I5").Value = 0.1760995803479087
$ws.Range("J5").Value = 0.1760995803479087
$ws.Range("M5").Value = 4.224096333333333
$ws.Range("N5").Value = 12.672289
$ws.Range("O5").Value = 0.1277189908446358
$ws.Range("P5").Value = 0.1277189908446358
$ws.Range("Q5").Value = 53.63457894833422
$ws.Range("R5").Value = 482.711210535008
$ws.Range("S5").Value = 0.02249126069019877
$ws.Range("T5").Value = 0.02249126069019877
$ws.Range("I6").Value = 0.1760995803479087
$ws.Range("J6").Value = 0.1760995803479087
$ws.Range("O6").Value = 0.4492078640046304
$ws.Range("P6").Value = 0.4492078640046304
$ws.Range("S6").Value = 0.07910531634019585
$ws.Range("T6").Value = 0.07910531634019587
$ws.Range("I7").Value = 0.1760995803479087
$ws.Range("J7").Value = 0.1760995803479087
$ws.Range("O7").Value = 0.4230731451507339
$ws.Range("P7").Value = 0.4230731451507338
$ws.Range("S7").Value = 0.0745030033175141
$ws.Range("T7").Value = 0.0745030033175141
$ws.Range("I8").Value = 0.5482939373535334
$ws.Range("J8").Value = 0.5482939373535334
$ws.Range("M8").Value = 4.224096333333333
$ws.Range("N8").Value = 12.672289
$ws.Range("O8").Value = 0.1277189908446358
$ws.Range("P8").Value = 0.1277189908446358
$ws.Range("Q8").Value = 166.9936658098932
$ws.Range("R8").Value = 1502.942992289039
$ws.Range("S8").Value = 0.07002754836502527
$ws.Range("T8").Value = 0.07002754836502527
$ws.Range("I9").Value = 0.5482939373535334
$ws.Range("J9").Value = 0.5482939373535334
$ws.Range("O9").Value = 0.4492078640046304
$ws.Range("P9").Value = 0.4492078640046304
$ws.Range("S9").Value = 0.2462979484452693
$ws.Range("T9").Value = 0.2462979484452693
$ws.Range("I10").Value = 0.5482939373535334
$ws.Range("J10").Value = 0.5482939373535334
$ws.Range("O10").Value = 0.4230731451507339
$ws.Range("P10").Value = 0.4230731451507338
$ws.Range("R10").Value = 4978.545590792102
$ws.Range("S10").Value = 0.2319684405432388
$ws.Range("T10").Value = 0.2319684405432388
